# =========================================================================
# AnalisisDatosI.xlsx - add "Estacionalidad" worksheet + forecasting rows
# =========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# -------------------------------------------------------------------
# 1. New worksheet "Estacionalidad", inserted right after sheet 1
# -------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Estacionalidad"

# -------------------------------------------------------------------
# 2. Text labels first (matches original authoring / shared-string
#    insertion order). D33/E33 use a leading quote so the "-"/"+"
#    prefixed text is stored as literal text (quotePrefix style).
# -------------------------------------------------------------------
$ws1.Range("B32").Value = "Crecimiento"
$ws1.Range("B33").Value = "Pronóstico con suavizado"
$ws1.Range("B31").Value = "Pronóstico lineal"
$ws1.Range("B34").Value = "Intervalo de confianza"
$ws1.Range("D33").Formula = "'-Intervalo"
$ws1.Range("E33").Formula = "'+Intervalo"

$ws2.Range("A21").Value = "Temporalidad"
$ws2.Range("A1").Value  = "Fechas"
$ws2.Range("B1").Value  = "Valores"

# -------------------------------------------------------------------
# 3. Sheet1 forecast block, rows 31-34 (values / formulas)
# -------------------------------------------------------------------
$ws1.Range("A31").Value = 2020

$ws1.Range("C31").Formula = "=FORECAST(A31,B8:B21,A8:A21)"
$ws1.Range("C31").NumberFormat = "#,##0"

$ws1.Range("C32").Formula = "=GROWTH(B8:B21,A8:A21,A31)"
$ws1.Range("C32").NumberFormat = "#,##0"

$ws1.Range("C33").Formula = "=_xlfn.FORECAST.ETS(A31,B8:B21,A8:A21)"
$ws1.Range("C33").NumberFormat = "#,##0"

$ws1.Range("C34").Formula = "=_xlfn.FORECAST.ETS.CONFINT(A31,B8:B21,A8:A21,0.9)"
$ws1.Range("C34").NumberFormat = "#,##0"

$ws1.Range("D34").Formula = "=C33-C34"
$ws1.Range("D34").NumberFormat = "#,##0"
$ws1.Range("E34").Formula = "=C33+C34"
$ws1.Range("E34").NumberFormat = "#,##0"

# -------------------------------------------------------------------
# 4. Sheet1 column I (rows 33-45) - left-over date-formatted column.
#    Build the style once, then clone it (format only) onto the rest
#    so every cell shares a single stylesheet entry.
# -------------------------------------------------------------------
$ws1.Range("I33").NumberFormat = "mm-dd-yy"
$ws1.Range("I33").Copy()
$ws1.Range("I34:I45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -------------------------------------------------------------------
# 5. Sheet2 ("Estacionalidad") - dates / values + ETS forecast table
# -------------------------------------------------------------------
$fechas  = @(42370,42371,42372,42373,42374,42375,42376,42377,42378,42379,42380,42381,42382)
$valores = @(5,7,5,12,13,6,4,8,11,13,6,5,7)

for ($i = 0; $i -lt $fechas.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $fechas[$i]
    $ws2.Cells.Item($row, 2).Value = $valores[$i]
}

# Date format for A2:A19 - create once on A2, clone the rest.
$ws2.Range("A2").NumberFormat = "mm-dd-yy"
$ws2.Range("A2").Copy()
$ws2.Range("A3:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Forecast rows 15-19 (dates continue, values predicted with FORECAST.ETS)
$ws2.Range("A15").Value = 42383
$ws2.Range("A16").Value = 42384
$ws2.Range("A17").Value = 42385
$ws2.Range("A18").Value = 42386
$ws2.Range("A19").Value = 42387

$ws2.Range("B15").Formula = "=_xlfn.FORECAST.ETS(A15,B2:B14,A2:A14)"
$ws2.Range("B16").Formula = "=_xlfn.FORECAST.ETS(A16,B`$2:B`$14,A`$2:A`$14)"
$ws2.Range("B17").Formula = "=_xlfn.FORECAST.ETS(A17,B`$2:B`$14,A`$2:A`$14)"
$ws2.Range("B18").Formula = "=_xlfn.FORECAST.ETS(A18,B`$2:B`$14,A`$2:A`$14)"
$ws2.Range("B19").Formula = "=_xlfn.FORECAST.ETS(A19,B`$2:B`$14,A`$2:A`$14)"

# 0.00 format for B15:B19 - create once, clone the rest.
$ws2.Range("B15").NumberFormat = "0.00"
$ws2.Range("B15").Copy()
$ws2.Range("B16:B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("B21").Formula = "=_xlfn.FORECAST.ETS.SEASONALITY(B2:B14,A2:A14)"

# -------------------------------------------------------------------
# 6. Column widths (sheet1 D & I, sheet2 A)
# -------------------------------------------------------------------
$ws1.Columns.Item(4).ColumnWidth = 11.86
$ws1.Columns.Item(9).ColumnWidth = 13.29
$ws2.Columns.Item(1).ColumnWidth = 13.29

# -------------------------------------------------------------------
# 7. Selections / active sheet (Estacionalidad ends up active, as in
#    the target workbook).
# -------------------------------------------------------------------
$ws1.Range("B35").Select()
$ws2.Range("A20").Select()
